$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.205.14"
$ws.Range("E2").Value = "  +0.86%  "

$ws.Range("D3").Value = "1.561.98"
$ws.Range("E3").Value = "  -0.15%  "

$ws.Range("E4").Value = "  -0.49%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.11"
$ws.Range("E5").Value = "  +1.06%  "

$ws.Range("E6").Value = "  -0.01%  "

$ws.Range("E7").Value = "  -0.57%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.01"
$ws.Range("E8").Value = "  -0.35%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.248"
$ws.Range("E9").Value = "  +0.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0596"
$ws.Range("E10").Value = "  -0.89%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0867"
$ws.Range("E11").Value = "  +1.15%  "

$ws.Range("D12").Value = "1.788.14"
$ws.Range("E12").Value = "  +0.08%  "

$ws.Range("D13").Value = "1.564.27"
$ws.Range("E13").Value = "  +0.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.77"
$ws.Range("E14").Value = "  +0.20%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.517"
$ws.Range("E15").Value = "  -0.59%  "

$ws.Range("D16").Value = "27.180.01"
$ws.Range("E16").Value = "  +0.78%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.89"
$ws.Range("E17").Value = "  -0.35%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.46"
$ws.Range("E18").Value = "  +1.54%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.60"
$ws.Range("E19").Value = "  -0.34%  "

$ws.Range("D20").Value = "0.0₃0701"
$ws.Range("E20").Value = "  -0.83%  "

$ws.Range("E21").Value = "  -0.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.13"
$ws.Range("E22").Value = "  +0.46%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.18"
$ws.Range("E23").Value = "  -0.53%  "

$ws.Range("E24").Value = "  -0.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.78"
$ws.Range("E25").Value = "  -0.30%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.62"
$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.00"
$ws.Range("E27").Value = "  -0.73%  "

$ws.Range("E28").Value = "  +1.13%  "

$ws.Range("E29").Value = "  -0.49%  "

$ws.Range("E30").Value = "  +1.82%  "

$ws.Range("E31").Value = "  -0.38%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.23"
$ws.Range("E32").Value = "  -0.11%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.15"
$ws.Range("E33").Value = "  +1.06%  "

$ws.Range("D34").Value = "1.438.32"
$ws.Range("E34").Value = "  +0.92%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.10"
$ws.Range("E35").Value = "  +3.55%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.60"
$ws.Range("E36").Value = "  -0.62%  "

$ws.Range("E37").Value = "  -0.02%  "

$ws.Range("E38").Value = "  +0.35%  "

$ws.Range("E39").Value = "  -0.51%  "

$ws.Range("E40").Value = "  +2.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.40"
$ws.Range("E41").Value = "  +2.90%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.806"
$ws.Range("E42").Value = "  -0.42%  "

$ws.Range("E43").Value = "  -0.41%  "

$ws.Range("E44").Value = "  -0.98%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.32"
$ws.Range("E45").Value = "  -0.90%  "

$ws.Range("E46").Value = "  -1.14%  "

$ws.Range("D47").Value = "1.704.31"
$ws.Range("E47").Value = "  +0.29%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.87"
$ws.Range("E48").Value = "  -1.83%  "

$ws.Range("E49").Value = "  +1.18%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0952"
$ws.Range("E50").Value = "  -0.73%  "

$ws.Range("D51").Value = "0.0₇0962"
$ws.Range("E51").Value = "  -5.59%  "
